# Generate Report for Handback
# Update the timestamp strings produced by the latest handback report run.

$wb = $excel.ActiveWorkbook

# --- "Overview" sheet: Latest HO Xliff Generate Date (col G) for first file ---
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-11-02 04:43:25"

# --- "zh-cn" sheet: Correspond Handoff Datetime (col H) / Handback DateTime (col K) ---
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-11-02 04:43:11"
$wsZhCn.Range("K2").Value = "2016-11-02 04:44:00"

# --- "de-de" sheet: Correspond Handoff Datetime (col H) shares the same
#     underlying timestamp text as Overview!G2, and Correspond Handback
#     DateTime (col K) gets its own new value ---
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-11-02 04:43:25"
$wsDeDe.Range("K2").Value = "2016-11-02 04:44:19"
